# live_trading_results.xlsx - apply trade #107 close / new-trade commit
#
# Helper: force a cell to be stored as a literal text string (number
# format "@") before assigning it, so that numeric-looking / date-looking
# / percent-looking strings (e.g. "69.5%", "2026-02-16", "21:42:11",
# "2.81") are NOT silently reinterpreted by Excel as a number/date/time.
function Set-TextCell($ws, $row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = [string]$val
}

# Helper: plain numeric cell.
function Set-NumCell($ws, $row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

# Helper: write one full "trade log" row (columns A..N) in one of the
# per-strategy sheets / the All Trades sheet.
function Set-TradeRow($ws, $row, $tradeNum, $date, $time, $strategy, $side, $entry, $exit, $status, $plPct, $plDollar, $conf, $entryReason, $exitReason, $durationMin) {
    Set-NumCell  $ws $row 1  $tradeNum
    Set-TextCell $ws $row 2  $date
    Set-TextCell $ws $row 3  $time
    Set-TextCell $ws $row 4  $strategy
    Set-TextCell $ws $row 5  $side
    Set-NumCell  $ws $row 6  $entry
    if ($null -ne $exit) {
        Set-NumCell $ws $row 7 $exit
    }
    Set-TextCell $ws $row 8  $status
    Set-NumCell  $ws $row 9  $plPct
    Set-NumCell  $ws $row 10 $plDollar
    Set-NumCell  $ws $row 11 $conf
    Set-TextCell $ws $row 12 $entryReason
    if ($null -ne $exitReason) {
        Set-TextCell $ws $row 13 $exitReason
    }
    Set-NumCell  $ws $row 14 $durationMin
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet: Summary
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

# Row 2 - OVERALL / ALL COMBINED
Set-NumCell  $wsSummary 2 3 82
Set-TextCell $wsSummary 2 4 "69.5%"
Set-TextCell $wsSummary 2 5 "+24.5193%"
Set-TextCell $wsSummary 2 6 "+0.2990%"

# Row 3 - STRATEGY / leadlag
Set-NumCell  $wsSummary 3 3 81
Set-TextCell $wsSummary 3 4 "46.9%"
Set-TextCell $wsSummary 3 5 "+13.1980%"
Set-TextCell $wsSummary 3 6 "+0.1629%"

# Row 4 - STRATEGY / momentum
Set-TextCell $wsSummary 4 4 "76.0%"
Set-TextCell $wsSummary 4 5 "+11.3213%"
Set-TextCell $wsSummary 4 6 "+0.4529%"

# ---------------------------------------------------------------------
# Sheet: leadlag
# ---------------------------------------------------------------------
$wsLeadlag = $wb.Worksheets.Item("leadlag")

# Close out trades #78-#82 (rows 58-62)
Set-NumCell  $wsLeadlag 58 7  68492.54893
Set-TextCell $wsLeadlag 58 8  "CLOSED"
Set-NumCell  $wsLeadlag 58 9  0.1788
Set-NumCell  $wsLeadlag 58 10 1.79
Set-TextCell $wsLeadlag 58 13 "time_exit_5min"
Set-NumCell  $wsLeadlag 58 14 5

Set-NumCell  $wsLeadlag 59 7  69393.66779399999
Set-TextCell $wsLeadlag 59 8  "CLOSED"
Set-NumCell  $wsLeadlag 59 9  -1.0782
Set-NumCell  $wsLeadlag 59 10 -10.78
Set-TextCell $wsLeadlag 59 13 "time_exit_5min"
Set-NumCell  $wsLeadlag 59 14 5

Set-NumCell  $wsLeadlag 60 7  69122.606841
Set-TextCell $wsLeadlag 60 8  "CLOSED"
Set-NumCell  $wsLeadlag 60 9  0.6687
Set-NumCell  $wsLeadlag 60 10 6.69
Set-TextCell $wsLeadlag 60 13 "time_exit_5min"
Set-NumCell  $wsLeadlag 60 14 5

Set-NumCell  $wsLeadlag 61 7  68623.57118499999
Set-TextCell $wsLeadlag 61 8  "CLOSED"
Set-NumCell  $wsLeadlag 61 9  -0.0448
Set-NumCell  $wsLeadlag 61 10 -0.45
Set-TextCell $wsLeadlag 61 13 "time_exit_5min"
Set-NumCell  $wsLeadlag 61 14 5

Set-NumCell  $wsLeadlag 62 7  68130.64425500001
Set-TextCell $wsLeadlag 62 8  "CLOSED"
Set-NumCell  $wsLeadlag 62 9  0.6593
Set-NumCell  $wsLeadlag 62 10 6.59
Set-TextCell $wsLeadlag 62 13 "time_exit_5min"
Set-NumCell  $wsLeadlag 62 14 5

# New trade #107 (row 83) - freshly opened
Set-TradeRow $wsLeadlag 83 107 "2026-02-16" "21:42:11" "leadlag" "UP" 68475.855 $null "OPEN" 0 0 0.7047 "Coinbase leading with 0.070% move" $null 0

# Column J (index 10) width 7 -> 8
$wsLeadlag.Columns.Item(10).ColumnWidth = 7.166666666666667

# ---------------------------------------------------------------------
# Sheet: momentum
# ---------------------------------------------------------------------
$wsMomentum = $wb.Worksheets.Item("momentum")

# Close out trades #76-#77 (rows 21-22)
Set-NumCell  $wsMomentum 21 7  69237.643658
Set-TextCell $wsMomentum 21 8  "CLOSED"
Set-NumCell  $wsMomentum 21 9  0.8126
Set-NumCell  $wsMomentum 21 10 8.130000000000001
Set-TextCell $wsMomentum 21 13 "time_exit_5min"
Set-NumCell  $wsMomentum 21 14 5

Set-NumCell  $wsMomentum 22 7  69005.10633
Set-TextCell $wsMomentum 22 8  "CLOSED"
Set-NumCell  $wsMomentum 22 9  0.4558
Set-NumCell  $wsMomentum 22 10 4.56
Set-TextCell $wsMomentum 22 13 "time_exit_5min"
Set-NumCell  $wsMomentum 22 14 5

# ---------------------------------------------------------------------
# Sheet: All Trades (rows 77-83 newly appended, mirroring the above
# closes plus the fresh trade's per-strategy row order)
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

Set-TradeRow $wsAll 77 76 "2026-02-16" "21:36:09" "momentum" "UP" 68679.56          69237.643658      "CLOSED" 0.8126  8.130000000000001 0.9    "Upward momentum: 0.287% over 10 samples" "time_exit_5min" 5
Set-TradeRow $wsAll 78 77 "2026-02-16" "21:36:16" "momentum" "UP" 68691.99000000001 69005.10633       "CLOSED" 0.4558  4.56               0.9    "Upward momentum: 0.197% over 10 samples" "time_exit_5min" 5
Set-TradeRow $wsAll 79 78 "2026-02-16" "21:36:27" "leadlag"  "DOWN" 68615.23         68492.54893       "CLOSED" 0.1788  1.79               0.75   "Binance leading with -0.148% move"       "time_exit_5min" 5
Set-TradeRow $wsAll 80 79 "2026-02-16" "21:36:33" "leadlag"  "DOWN" 68653.47500000001 69393.66779399999 "CLOSED" -1.0782 -10.78             0.6561 "Coinbase leading with -0.066% move"      "time_exit_5min" 5
Set-TradeRow $wsAll 81 80 "2026-02-16" "21:36:40" "leadlag"  "UP"   68663.435        69122.606841      "CLOSED" 0.6687  6.69               0.75   "Binance leading with 0.097% move"        "time_exit_5min" 5
Set-TradeRow $wsAll 82 81 "2026-02-16" "21:36:46" "leadlag"  "DOWN" 68592.83         68623.57118499999 "CLOSED" -0.0448 -0.45              0.75   "Binance leading with -0.097% move"       "time_exit_5min" 5
Set-TradeRow $wsAll 83 82 "2026-02-16" "21:36:52" "leadlag"  "DOWN" 68582.78999999999 68130.64425500001 "CLOSED" 0.6593  6.59               0.75   "Binance leading with -0.118% move"       "time_exit_5min" 5

# Column J (index 10) width 7 -> 8
$wsAll.Columns.Item(10).ColumnWidth = 7.166666666666667

# ---------------------------------------------------------------------
# Sheet: Comparison
# ---------------------------------------------------------------------
$wsCompare = $wb.Worksheets.Item("Comparison")

# Row 2 - leadlag
Set-NumCell  $wsCompare 2 2 81
Set-TextCell $wsCompare 2 3 "46.9%"
Set-TextCell $wsCompare 2 4 "2.81"
Set-TextCell $wsCompare 2 5 "+0.5390%"
Set-TextCell $wsCompare 2 6 "-0.3167%"
Set-TextCell $wsCompare 2 7 "1.70"
Set-TextCell $wsCompare 2 8 "-1.0782%"

# Row 3 - momentum
Set-TextCell $wsCompare 3 3 "76.0%"
Set-TextCell $wsCompare 3 4 "11.07"
Set-TextCell $wsCompare 3 5 "+0.6550%"
